$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4339.821745983342
$ws.Range("C3").Value = 4210.617293134998
$ws.Range("C4").Value = 4141.63110474518
$ws.Range("C5").Value = 4141.63110474518
$ws.Range("C6").Value = 4141.63110474518
$ws.Range("C7").Value = 4141.63110474518
$ws.Range("C8").Value = 4141.63110474518
$ws.Range("C9").Value = 4118.3382300992
$ws.Range("C10").Value = 4118.3382300992
$ws.Range("C11").Value = 4118.3382300992
$ws.Range("C12").Value = 4027.95152112983
